$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 12 with value continuing the A column time-log series,
# matching the number format/style used by the rest of the A column.
$ws.Cells.Item(12, 1).Value = 0.027164351851851853
$ws.Range("A12").NumberFormat = $ws.Range("A11").NumberFormat

# Update the SUM formulas so they include the new row 12
$ws.Range("C2").Formula = "=SUM(A2:A12)"
$ws.Range("B3").Formula = "=SUM(A9:A12)"

# Update the selected / active cell in the sheet view
$ws.Range("D3").Select()

# Update the workbook window position (yWindow) to match the saved view
$wb.Windows.Item(1).Top = 456
